$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DOMA-3439: fix typo in Properties Export Template header
# "{d.i18n.tocketsInWork}" -> "{d.i18n.ticketsInWork}"
$ws.Range("E1").Value = "{d.i18n.ticketsInWork}"
